$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 1. Collapse the three verbose "CORE COMPETENCIES" paragraphs into a single
#    short summary line.
# ---------------------------------------------------------------------------
$coreCompetencies = $d.Paragraphs.Item(6)
$coreCompetencies.Range.Text = "Product Marketing Core " + $bullet + " Research & Analytics " + $bullet + " Communication & Technology"

# Remove the two paragraphs that followed it (now both sitting at index 7,
# since deleting the first collapses the next one into its place).
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# ---------------------------------------------------------------------------
# 2. Add a new "TECHNICAL SKILLS" section (heading + 3 detail paragraphs)
#    right after the "Managed national polling team ..." bullet, before the
#    closing LinkedIn/site sentence.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(47)
if ($anchor.Range.Text -notmatch "Managed national polling team") {
    throw "Anchor paragraph mismatch: " + $anchor.Range.Text
}

# Insert four fresh paragraphs after the anchor. Since the anchor is a plain
# (Normal-style) paragraph, every one of these inherits Normal style, so no
# spurious style leaks into the body paragraphs.
$anchor.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Item(48)

$heading.Range.InsertParagraphAfter()
$para1 = $d.Paragraphs.Item(49)

$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item(50)

$para2.Range.InsertParagraphAfter()
$para3 = $d.Paragraphs.Item(51)

$para1.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development"
$para2.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; A/B Testing & Conversion Optimization"
$para3.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Client Relationship Management & Business Development"

# Set heading text/style last, so its Heading2 formatting has nothing left
# to leak forward into.
$heading.Range.Text = "TECHNICAL SKILLS"
$heading.Style = "Heading 2"

Write-Output "Applied core-competencies collapse and TECHNICAL SKILLS section insert."
